$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 495, shifting existing rows 495-538 down to 496-539.
$ws.Range("A495").EntireRow.Insert()

# Populate the newly inserted row 495 with the new data record.
$ws.Range("A495").Value = 3
$ws.Range("B495").Value = "Femacal de La Calera"
$ws.Range("C495").Value = "Coquimbo"
$ws.Range("D495").Value = 45223
$ws.Range("E495").Value = 5
$ws.Range("F495").Value = 100112001
$ws.Range("G495").Value = "Berenjena"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 100
$ws.Range("K495").Value = 8500
$ws.Range("L495").Value = 9000
$ws.Range("M495").Value = 8750
$ws.Range("N495").Value = "`$/caja 60 unidades"
$ws.Range("O495").Value = "Región de Arica y Parinacota"
$ws.Range("P495").Value = 146
$ws.Range("Q495").Value = 60
$ws.Range("R495").Value = "Hortaliza"
